$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$tbl = $ws.ListObjects.Item(1)

# ---------------------------------------------------------------------------
# 1) Insert a new row for "Elixinol" right before the current row 83
#    (Zelda Therapeutics), i.e. physically insert a worksheet row at 83 and
#    push everything below it down by one.
# ---------------------------------------------------------------------------
$ws.Rows(83).Insert()

# Copy formatting from the similar existing "Cannabis US" row (row 81, which
# is not affected by the insert) onto the newly created row 83.
$ws.Range("A81:I81").Copy()
$ws.Range("A83:I83").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(83, 1).Value = 82
$ws.Cells.Item(83, 2).Value = "Cannabis US"
$ws.Cells.Item(83, 3).Value = "Elixinol"
$ws.Cells.Item(83, 4).Value = "Elixinol Global Ltd"
$ws.Cells.Item(83, 6).Value = "ELLXF"

# ---------------------------------------------------------------------------
# 2) Insert a new row for "Plus Product" right before what is now row 95
#    (originally row 94: "Cannabis AU" / "MMJ Group"), after the first
#    insertion shifted everything down by one.
# ---------------------------------------------------------------------------
$ws.Rows(95).Insert()

$ws.Range("A84:I84").Copy()
$ws.Range("A95:I95").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(95, 1).Value = 94
$ws.Cells.Item(95, 2).Value = "Cannabis US"
$ws.Cells.Item(95, 3).Value = "Plus Product"
$ws.Cells.Item(95, 4).Value = "Plus Product Inc"
$ws.Cells.Item(95, 6).Value = "PLPRF"

# ---------------------------------------------------------------------------
# 3) Renumber the "Nr." column so it stays sequential (1 .. 96) after the two
#    row insertions shifted the existing numbering down.
# ---------------------------------------------------------------------------
for ($r = 84; $r -le 97; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 1
}

# ---------------------------------------------------------------------------
# 4) Make sure the table (ListObject) covers the two new rows.
# ---------------------------------------------------------------------------
$tbl.Resize($ws.Range("A1:I97"))

# ---------------------------------------------------------------------------
# 5) Restore the selection shown in the bottom-right frozen pane.
# ---------------------------------------------------------------------------
$ws.Range("B86").Select()
